$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37/38 swap: Dai moves up to rank 35 (row37), Fetch.AI moves to rank 36 (row38) ---
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.17%  "

$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.16%  "

# --- Price/volume updates for remaining rows ---
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.625.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.497.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.26%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.41%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.199"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.648"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.05%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.44%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000299"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.81%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.88%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.076.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "592.96"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.17%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.814.73"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "12.73"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.96%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.507.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.66%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.983"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.76%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.55%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.79%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.95%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.49%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.91%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.50"
$ws.Range("D28").Style = "Normal"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.38%  "

# Row 33
$ws.Range("E33").Value = "  -0.28%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.712.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.59%  "

# Row 39
$ws.Range("E39").Value = "  -0.80%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.389"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.00%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "483.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.30%  "

# Row 43
$ws.Range("E43").Value = "  -3.75%  "

# Row 44
$ws.Range("E44").Value = "  -1.71%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.139"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.61%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "

# Row 48
$ws.Range("E48").Value = "  +0.34%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.33%  "

# Row 50
$ws.Range("E50").Value = "  +0.93%  "

# Row 51
$ws.Range("E51").Value = "  +10.59%  "
